# Auto-generated Excel COM-interop script applying the numeric value
# updates described by the commit diff for Sheets/Moogle_Profits.xlsx.
# Each ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheet gets its H:N market-data
# columns (currentAveragePrice*, LevePrice*, LeveProfit*) refreshed to the
# values captured by the scheduled Universalis price-sync runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 678274.3
$ws.Range("I15").Value = 678274.3
$ws.Range("K15").Value = 2034822.9
$ws.Range("M15").Value = -2034653.9
$ws.Range("H17").Value = 2164.5625
$ws.Range("J17").Value = 2164.5625
$ws.Range("L17").Value = 6493.6875
$ws.Range("N17").Value = -6829.6875
$ws.Range("H19").Value = 1112.5
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 1300
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 1300
$ws.Range("M19").Value = -825
$ws.Range("N19").Value = -1650
$ws.Range("H40").Value = 2192.84
$ws.Range("I40").Value = 1461.625
$ws.Range("J40").Value = 3492.7778
$ws.Range("K40").Value = 1461.625
$ws.Range("L40").Value = 3492.7778
$ws.Range("M40").Value = -1286.625
$ws.Range("N40").Value = -3842.7778
$ws.Range("H74").Value = 10105.529
$ws.Range("I74").Value = 7985.357
$ws.Range("K74").Value = 7985.357
$ws.Range("M74").Value = -7049.357
$ws.Range("H77").Value = 10105.529
$ws.Range("I77").Value = 7985.357
$ws.Range("K77").Value = 39926.785
$ws.Range("M77").Value = -35246.785
$ws.Range("H113").Value = 4959.4546
$ws.Range("I113").Value = 3811.2
$ws.Range("K113").Value = 3811.2
$ws.Range("M113").Value = -557.1999999999998
$ws.Range("H132").Value = 2149.8838
$ws.Range("I132").Value = 1908.9487
$ws.Range("K132").Value = 5726.8461
$ws.Range("M132").Value = -3196.8461
$ws.Range("H137").Value = 1936.0408
$ws.Range("I137").Value = 2003.7317
$ws.Range("J137").Value = 1589.125
$ws.Range("K137").Value = 6011.1951
$ws.Range("L137").Value = 4767.375
$ws.Range("M137").Value = -3461.1951
$ws.Range("N137").Value = -9867.375
$ws.Range("H138").Value = 7309.3037
$ws.Range("I138").Value = 5477.952
$ws.Range("K138").Value = 16433.856
$ws.Range("M138").Value = -11293.856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 59.6
$ws.Range("I4").Value = 59.6
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 59.6
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 56.4
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 331.33334
$ws.Range("I5").Value = 320.17648
$ws.Range("K5").Value = 320.17648
$ws.Range("M5").Value = -208.17648
$ws.Range("H32").Value = 2641.03
$ws.Range("I32").Value = 1852.2446
$ws.Range("J32").Value = 14998.667
$ws.Range("K32").Value = 1852.2446
$ws.Range("L32").Value = 14998.667
$ws.Range("M32").Value = -1565.2446
$ws.Range("N32").Value = -15572.667
$ws.Range("H45").Value = 3399
$ws.Range("I45").Value = 1748.875
$ws.Range("K45").Value = 1748.875
$ws.Range("M45").Value = -1371.875
$ws.Range("H74").Value = 9437303
$ws.Range("I74").Value = 10872130
$ws.Range("K74").Value = 10872130
$ws.Range("M74").Value = -10871256
$ws.Range("H77").Value = 9437303
$ws.Range("I77").Value = 10872130
$ws.Range("K77").Value = 54360650
$ws.Range("M77").Value = -54356282
$ws.Range("H102").Value = 2405.2307
$ws.Range("I102").Value = 2405.2307
$ws.Range("K102").Value = 2405.2307
$ws.Range("M102").Value = -783.2307000000001
$ws.Range("H110").Value = 3361.6875
$ws.Range("I110").Value = 3499.0715
$ws.Range("J110").Value = 2400
$ws.Range("K110").Value = 3499.0715
$ws.Range("L110").Value = 2400
$ws.Range("M110").Value = -1454.0715
$ws.Range("N110").Value = -6490
$ws.Range("H132").Value = 3567.681
$ws.Range("I132").Value = 2175.2104
$ws.Range("J132").Value = 9447
$ws.Range("K132").Value = 6525.6312
$ws.Range("L132").Value = 28341
$ws.Range("M132").Value = -3995.6312
$ws.Range("N132").Value = -33401

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 331.33334
$ws.Range("I4").Value = 320.17648
$ws.Range("K4").Value = 320.17648
$ws.Range("M4").Value = -205.17648
$ws.Range("H99").Value = 1506.6666
$ws.Range("I99").Value = 1200.5625
$ws.Range("J99").Value = 3955.5
$ws.Range("K99").Value = 1200.5625
$ws.Range("L99").Value = 3955.5
$ws.Range("M99").Value = 297.4375
$ws.Range("N99").Value = -6951.5
$ws.Range("H105").Value = 5082.6665
$ws.Range("I105").Value = 5097
$ws.Range("J105").Value = 5011
$ws.Range("K105").Value = 5097
$ws.Range("L105").Value = 5011
$ws.Range("M105").Value = -3350
$ws.Range("N105").Value = -8505
$ws.Range("H107").Value = 7333.3335
$ws.Range("J107").Value = 10000
$ws.Range("L107").Value = 10000
$ws.Range("N107").Value = -13840
$ws.Range("H129").Value = 92000
$ws.Range("J129").Value = 92000
$ws.Range("L129").Value = 92000
$ws.Range("N129").Value = -102000
$ws.Range("H134").Value = 2309.16
$ws.Range("I134").Value = 1655.375
$ws.Range("K134").Value = 4966.125
$ws.Range("M134").Value = -2431.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6749.2285
$ws.Range("I31").Value = 2761.7917
$ws.Range("J31").Value = 15449.091
$ws.Range("K31").Value = 2761.7917
$ws.Range("L31").Value = 15449.091
$ws.Range("M31").Value = -2466.7917
$ws.Range("N31").Value = -16039.091
$ws.Range("H34").Value = 6749.2285
$ws.Range("I34").Value = 2761.7917
$ws.Range("J34").Value = 15449.091
$ws.Range("K34").Value = 2761.7917
$ws.Range("L34").Value = 15449.091
$ws.Range("M34").Value = -2559.7917
$ws.Range("N34").Value = -15853.091
$ws.Range("H94").Value = 6464.8
$ws.Range("I94").Value = 6465.6
$ws.Range("J94").Value = 6464
$ws.Range("K94").Value = 6465.6
$ws.Range("L94").Value = 6464
$ws.Range("M94").Value = -6014.6
$ws.Range("N94").Value = -7366
$ws.Range("H105").Value = 1844.6666
$ws.Range("I105").Value = 1943.1428
$ws.Range("K105").Value = 1943.1428
$ws.Range("M105").Value = -196.1428000000001
$ws.Range("H122").Value = 1608.5714
$ws.Range("I122").Value = 1650.9
$ws.Range("J122").Value = 1502.75
$ws.Range("K122").Value = 4952.700000000001
$ws.Range("L122").Value = 4508.25
$ws.Range("M122").Value = -2502.700000000001
$ws.Range("N122").Value = -9408.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 34433.168
$ws.Range("I116").Value = 1349.75
$ws.Range("K116").Value = 4049.25
$ws.Range("M116").Value = -607.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7057.116
$ws.Range("I70").Value = 7500
$ws.Range("J70").Value = 7046.5713
$ws.Range("K70").Value = 7500
$ws.Range("L70").Value = 7046.5713
$ws.Range("M70").Value = -7230
$ws.Range("N70").Value = -7586.5713
$ws.Range("H73").Value = 7057.116
$ws.Range("I73").Value = 7500
$ws.Range("J73").Value = 7046.5713
$ws.Range("K73").Value = 7500
$ws.Range("L73").Value = 7046.5713
$ws.Range("M73").Value = -6564
$ws.Range("N73").Value = -8918.5713
$ws.Range("H113").Value = 3996.5
$ws.Range("I113").Value = 994.75
$ws.Range("K113").Value = 994.75
$ws.Range("M113").Value = 1175.25
$ws.Range("H122").Value = 10784.5
$ws.Range("I122").Value = 8499.5
$ws.Range("J122").Value = 11927
$ws.Range("K122").Value = 25498.5
$ws.Range("L122").Value = 35781
$ws.Range("M122").Value = -23048.5
$ws.Range("N122").Value = -40681
$ws.Range("H127").Value = 27217
$ws.Range("J127").Value = 27217
$ws.Range("L127").Value = 27217
$ws.Range("N127").Value = -37137
$ws.Range("H132").Value = 4103.7314
$ws.Range("I132").Value = 3657.426
$ws.Range("J132").Value = 5957.615
$ws.Range("K132").Value = 10972.278
$ws.Range("L132").Value = 17872.845
$ws.Range("M132").Value = -8442.278
$ws.Range("N132").Value = -22932.845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1604.8572
$ws.Range("J22").Value = 2547
$ws.Range("L22").Value = 2547
$ws.Range("N22").Value = -3137
$ws.Range("H27").Value = 1604.8572
$ws.Range("J27").Value = 2547
$ws.Range("L27").Value = 2547
$ws.Range("N27").Value = -2761
$ws.Range("H46").Value = 3764.6
$ws.Range("I46").Value = 947.6
$ws.Range("J46").Value = 5173.1
$ws.Range("K46").Value = 947.6
$ws.Range("L46").Value = 5173.1
$ws.Range("M46").Value = -759.6
$ws.Range("N46").Value = -5549.1
$ws.Range("H93").Value = 1710.7273
$ws.Range("I93").Value = 1631.8
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 1631.8
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -383.8
$ws.Range("N93").Value = -4996
$ws.Range("H100").Value = 2100.8
$ws.Range("I100").Value = 1929.625
$ws.Range("J100").Value = 2296.4285
$ws.Range("K100").Value = 1929.625
$ws.Range("L100").Value = 2296.4285
$ws.Range("M100").Value = -1388.625
$ws.Range("N100").Value = -3378.4285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 19000
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 19000
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
$ws.Range("H45").Value = 17183.334
$ws.Range("I45").Value = 30569
$ws.Range("J45").Value = 10490.5
$ws.Range("K45").Value = 30569
$ws.Range("L45").Value = 10490.5
$ws.Range("M45").Value = -30078
$ws.Range("N45").Value = -11472.5
$ws.Range("H54").Value = 61877.168
$ws.Range("J54").Value = 82088.336
$ws.Range("L54").Value = 82088.336
$ws.Range("N54").Value = -83128.336
$ws.Range("H125").Value = 70333
$ws.Range("J125").Value = 70333
$ws.Range("L125").Value = 70333
$ws.Range("N125").Value = -80173
$ws.Range("H132").Value = 5868.8823
$ws.Range("I132").Value = 3125.5
$ws.Range("J132").Value = 8307.444
$ws.Range("K132").Value = 9376.5
$ws.Range("L132").Value = 24922.332
$ws.Range("M132").Value = -6846.5
$ws.Range("N132").Value = -29982.332
$ws.Range("H136").Value = 4675.875
$ws.Range("I136").Value = 3772.0715
$ws.Range("K136").Value = 11316.2145
$ws.Range("M136").Value = -8766.2145
